# The sheet holds weekly price records for "Arveja Verde" at Terminal
# Hortofrutícola Agro Chillán. A new weekly record needs to be inserted
# as row 18 (pushing the existing row 18 and all following rows down by
# one), growing the used range from A1:R44 to A1:R45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18; everything from the old row 18
# onward shifts down to make room (old row 18 -> new row 19, ..., old
# row 44 -> new row 45).
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C18").Value = "Ñuble"
$ws.Range("D18").Value = 44557
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 100112022
$ws.Range("G18").Value = "Arveja Verde"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 17000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 17500
$ws.Range("N18").Value = "$/saco 25 kilos"
$ws.Range("O18").Value = "Provincia de Diguillín"
$ws.Range("P18").Value = 700
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of the sheet
# (yyyy-mm-dd hh:mm:ss style used for every other row's date cell).
$ws.Range("D18").NumberFormat = $ws.Range("D19").NumberFormat()
